$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1684587813620072
$ws.Range("C2").Value = 0.6057347670250897
$ws.Range("J2").Value = 0.01433691756272401
$ws.Range("P2").Value = 0.1433691756272401
$ws.Range("S2").Value = 0.06810035842293907
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.03910614525139665
$ws.Range("J3").Value = 0.0111731843575419
$ws.Range("P3").Value = 0.7486033519553073
$ws.Range("S3").Value = 0.1955307262569832
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("O4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.7111111111111111
$ws.Range("S4").Value = 0.2222222222222222
$ws.Range("B6").Value = 0.06779661016949153
$ws.Range("D6").Value = 0.008474576271186441
$ws.Range("F6").Value = 0.1059322033898305
$ws.Range("J6").Value = 0.2245762711864407
$ws.Range("O6").Value = 0.02542372881355932
$ws.Range("Q6").Value = 0.1355932203389831
$ws.Range("R6").Value = 0.05932203389830509
$ws.Range("S6").Value = 0.3728813559322034
$ws.Range("B7").Value = 0.1075268817204301
$ws.Range("D7").Value = 0.03225806451612903
$ws.Range("E7").Value = 0.01075268817204301
$ws.Range("F7").Value = 0.09677419354838709
$ws.Range("J7").Value = 0.1075268817204301
$ws.Range("O7").Value = 0.01075268817204301
$ws.Range("Q7").Value = 0.1021505376344086
$ws.Range("R7").Value = 0.1021505376344086
$ws.Range("S7").Value = 0.4301075268817204
$ws.Range("B8").Value = 0.1170731707317073
$ws.Range("D8").Value = 0.01463414634146342
$ws.Range("F8").Value = 0.08292682926829269
$ws.Range("J8").Value = 0.0975609756097561
$ws.Range("O8").Value = 0.01951219512195122
$ws.Range("Q8").Value = 0.1853658536585366
$ws.Range("R8").Value = 0.0975609756097561
$ws.Range("S8").Value = 0.3853658536585366
$ws.Range("B9").Value = 0.1222707423580786
$ws.Range("D9").Value = 0.008733624454148471
$ws.Range("E9").Value = 0.004366812227074236
$ws.Range("F9").Value = 0.05240174672489083
$ws.Range("J9").Value = 0.1004366812227074
$ws.Range("O9").Value = 0.02620087336244541
$ws.Range("Q9").Value = 0.1746724890829694
$ws.Range("R9").Value = 0.08733624454148471
$ws.Range("S9").Value = 0.4235807860262009
$ws.Range("B10").Value = 0.09815436241610738
$ws.Range("D10").Value = 0.02768456375838926
$ws.Range("E10").Value = 0.005033557046979865
$ws.Range("F10").Value = 0.07130872483221476
$ws.Range("J10").Value = 0.1107382550335571
$ws.Range("O10").Value = 0.01845637583892618
$ws.Range("Q10").Value = 0.2315436241610738
$ws.Range("R10").Value = 0.08305369127516779
$ws.Range("S10").Value = 0.3540268456375839
$ws.Range("G11").Value = 0.1208791208791209
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.1611721611721612
$ws.Range("L11").Value = 0.6336996336996337
$ws.Range("S11").Value = 0.007326007326007326
$ws.Range("G12").Value = 0.7231638418079096
$ws.Range("J12").Value = 0.1581920903954802
$ws.Range("K12").Value = 0.01694915254237288
$ws.Range("L12").Value = 0.03389830508474576
$ws.Range("S12").Value = 0.06779661016949153
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.004329004329004329
$ws.Range("H15").Value = 0.1298701298701299
$ws.Range("I15").Value = 0.06493506493506493
$ws.Range("J15").Value = 0.3463203463203463
$ws.Range("K15").Value = 0.0735930735930736
$ws.Range("M15").Value = 0.01298701298701299
$ws.Range("O15").Value = 0.1082251082251082
$ws.Range("S15").Value = 0.2597402597402597
$ws.Range("F16").Value = 0.004950495049504951
$ws.Range("H16").Value = 0.1584158415841584
$ws.Range("I16").Value = 0.09900990099009901
$ws.Range("J16").Value = 0.4455445544554456
$ws.Range("K16").Value = 0.1188118811881188
$ws.Range("M16").Value = 0.01485148514851485
$ws.Range("N16").Value = 0.009900990099009901
$ws.Range("O16").Value = 0.05445544554455446
$ws.Range("S16").Value = 0.09405940594059406
$ws.Range("F17").Value = 0.01363636363636364
$ws.Range("H17").Value = 0.1636363636363636
$ws.Range("I17").Value = 0.09545454545454546
$ws.Range("J17").Value = 0.3886363636363636
$ws.Range("K17").Value = 0.09318181818181819
$ws.Range("M17").Value = 0.02727272727272727
$ws.Range("O17").Value = 0.09090909090909091
$ws.Range("S17").Value = 0.1272727272727273
$ws.Range("F18").Value = 0.03174603174603174
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.1164021164021164
$ws.Range("J18").Value = 0.4021164021164021
$ws.Range("K18").Value = 0.08994708994708994
$ws.Range("M18").Value = 0.02116402116402116
$ws.Range("O18").Value = 0.07407407407407407
$ws.Range("S18").Value = 0.1216931216931217
$ws.Range("F19").Value = 0.01714285714285714
$ws.Range("H19").Value = 0.2040816326530612
$ws.Range("I19").Value = 0.1085714285714286
$ws.Range("J19").Value = 0.3657142857142857
$ws.Range("K19").Value = 0.09877551020408164
$ws.Range("M19").Value = 0.0236734693877551
$ws.Range("N19").Value = 0.0008163265306122449
$ws.Range("O19").Value = 0.0563265306122449
$ws.Range("S19").Value = 0.1248979591836735
